# Update Data Dictionary mit Temperatur
# Adds a new "Temperatur" (temperature) section to Tabelle1, mirroring the
# existing "Bewölkung" section, and fixes the header-row styling of the
# "Bewölkung" block (B19/C19) to match the bold+centred header style used
# elsewhere (e.g. B10/C10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

function Set-FormatLike {
    param($targetCell, $sourceCell)
    $ws.Range($sourceCell).Copy() | Out-Null
    $ws.Range($targetCell).PasteSpecial(-4122) | Out-Null
}

# --- Fix existing "Bewölkung" header row (B19/C19): center+bold like B10/C10 ---
Set-FormatLike "B19" "B10"
Set-FormatLike "C19" "B10"

# --- New "Temperatur" section (rows 25-31), mirroring "Bewölkung" (rows 18-23) ---

# Section title
$ws.Range("A25").Value = "Temperatur"
Set-FormatLike "A25" "A18"

# Column header row
$ws.Range("A26").Value = "Bezeichnung"
Set-FormatLike "A26" "A18"
$ws.Range("A26").Font.Bold = $true

$ws.Range("B26").Value = "Temperatur"
Set-FormatLike "B26" "A2"

$ws.Range("C26").Value = "Eigener Code"
Set-FormatLike "C26" "A2"

# Data rows
$ws.Range("A27").Value = "Eistag"
Set-FormatLike "A27" "A18"
$ws.Range("B27").Value = "-10 - 0"
Set-FormatLike "B27" "C4"
$ws.Range("C27").Value = 0
Set-FormatLike "C27" "C3"

$ws.Range("A28").Value = "Kalttag"
Set-FormatLike "A28" "A18"
$ws.Range("C28").Value = 1
Set-FormatLike "C28" "C3"

$ws.Range("A29").Value = "Warmtag"
Set-FormatLike "A29" "A18"
$ws.Range("B29").Value = "16-24"
Set-FormatLike "B29" "C4"
$ws.Range("C29").Value = 2
Set-FormatLike "C29" "C3"

$ws.Range("A30").Value = "Sommertag"
Set-FormatLike "A30" "A18"
$ws.Range("B30").Value = "25-30"
Set-FormatLike "B30" "C4"
$ws.Range("C30").Value = 3
Set-FormatLike "C30" "C3"

$ws.Range("A31").Value = "Hitzetag"
Set-FormatLike "A31" "A18"
$ws.Range("B31").Value = "31-40"
Set-FormatLike "B31" "C4"
$ws.Range("C31").Value = 4
Set-FormatLike "C31" "C3"

# B28 filled in last (matches the original author's edit order)
$ws.Range("B28").Value = "1-15"
Set-FormatLike "B28" "C4"

$excel.CutCopyMode = $false

# --- Update the sheet view to reflect where the author ended up working ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D20").Select()
